# Apply the recorded edit to 建筑业企业资产负债情况.xlsx
#
# Summary of the change:
#  1. Columns H:M (the non-cumulative duplicate series) are removed entirely,
#     shrinking the used range from A1:M16 to A1:G16.
#  2. Within columns A:G, the data rows for the "B" and "C" sub-series of each
#     year are swapped in place:
#       row 3  <-> row 4   (2012年B <-> 2012年C)
#       row 6  <-> row 7   (2013年B <-> 2013年C)
#       row 9  <-> row 10  (2014年B <-> 2014年C)
#       row 12 <-> row 13  (2015年B <-> 2015年C)
#     All other rows (1,2,5,8,11,14,15,16) stay untouched.
#
# Column E is blank (no value) for rows 9,10,12,13 on both sides of the swap,
# so it is handled separately and only touched where it actually holds data,
# to avoid needlessly clearing/recreating empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove columns H:M (only A:G should remain) ---
$ws.Range("H1:M16").Delete()

# --- Step 2: swap the row content for each paired year ---
function Swap-Rows($sheet, $r1, $r2) {
    # Columns A:D
    $adRange1 = "A" + $r1 + ":D" + $r1
    $adRange2 = "A" + $r2 + ":D" + $r2
    $ad1 = $sheet.Range($adRange1).Value2
    $ad2 = $sheet.Range($adRange2).Value2
    $sheet.Range($adRange1).Value2 = $ad2
    $sheet.Range($adRange2).Value2 = $ad1

    # Column E (only swap if at least one side actually has a value)
    $eCell1 = $sheet.Range("E" + $r1)
    $eCell2 = $sheet.Range("E" + $r2)
    $e1 = $eCell1.Value2
    $e2 = $eCell2.Value2
    if (($e1 -ne "") -or ($e2 -ne "")) {
        $eCell1.Value2 = $e2
        $eCell2.Value2 = $e1
    }

    # Columns F:G
    $fgRange1 = "F" + $r1 + ":G" + $r1
    $fgRange2 = "F" + $r2 + ":G" + $r2
    $fg1 = $sheet.Range($fgRange1).Value2
    $fg2 = $sheet.Range($fgRange2).Value2
    $sheet.Range($fgRange1).Value2 = $fg2
    $sheet.Range($fgRange2).Value2 = $fg1
}

Swap-Rows $ws 3 4
Swap-Rows $ws 6 7
Swap-Rows $ws 9 10
Swap-Rows $ws 12 13
